$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1418
$ws.Cells.Item(1418, 1).Value = 44448
$ws.Cells.Item(1418, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1418, 2).Value = "'9061618"
$ws.Cells.Item(1418, 3).Value = 3011
$ws.Cells.Item(1418, 4).Value = "Order 9061618 Card(Stripe)"
$ws.Cells.Item(1418, 6).Value = 883.93

# Row 1419
$ws.Cells.Item(1419, 1).Value = 44448
$ws.Cells.Item(1419, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1419, 2).Value = "'9061618"
$ws.Cells.Item(1419, 3).Value = 2611
$ws.Cells.Item(1419, 4).Value = "Order 9061618 Card(Stripe)"
$ws.Cells.Item(1419, 6).Value = 106.07

# Row 1420
$ws.Cells.Item(1420, 1).Value = 44448
$ws.Cells.Item(1420, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1420, 2).Value = "'9061618"
$ws.Cells.Item(1420, 3).Value = 1930
$ws.Cells.Item(1420, 4).Value = "Order 9061618 Card(Stripe)"
$ws.Cells.Item(1420, 5).Value = 990

# Row 1421
$ws.Cells.Item(1421, 1).Value = 44448
$ws.Cells.Item(1421, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1421, 2).Value = "'7061821"
$ws.Cells.Item(1421, 3).Value = 3011
$ws.Cells.Item(1421, 4).Value = "Order 7061821 Card(Stripe)"
$ws.Cells.Item(1421, 6).Value = 1132.14

# Row 1422
$ws.Cells.Item(1422, 1).Value = 44448
$ws.Cells.Item(1422, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1422, 2).Value = "'7061821"
$ws.Cells.Item(1422, 3).Value = 2611
$ws.Cells.Item(1422, 4).Value = "Order 7061821 Card(Stripe)"
$ws.Cells.Item(1422, 6).Value = 135.86

# Row 1423
$ws.Cells.Item(1423, 1).Value = 44448
$ws.Cells.Item(1423, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1423, 2).Value = "'7061821"
$ws.Cells.Item(1423, 3).Value = 1930
$ws.Cells.Item(1423, 4).Value = "Order 7061821 Card(Stripe)"
$ws.Cells.Item(1423, 5).Value = 1268

# Row 1424
$ws.Cells.Item(1424, 1).Value = 44448
$ws.Cells.Item(1424, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1424, 2).Value = "'3092057"
$ws.Cells.Item(1424, 3).Value = 3011
$ws.Cells.Item(1424, 4).Value = "Order 3092057 Swish +46703564388"
$ws.Cells.Item(1424, 6).Value = 1062.5

# Row 1425
$ws.Cells.Item(1425, 1).Value = 44448
$ws.Cells.Item(1425, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1425, 2).Value = "'3092057"
$ws.Cells.Item(1425, 3).Value = 2611
$ws.Cells.Item(1425, 4).Value = "Order 3092057 Swish +46703564388"
$ws.Cells.Item(1425, 6).Value = 127.5

# Row 1426
$ws.Cells.Item(1426, 1).Value = 44448
$ws.Cells.Item(1426, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1426, 2).Value = "'3092057"
$ws.Cells.Item(1426, 3).Value = 1930
$ws.Cells.Item(1426, 4).Value = "Order 3092057 Swish +46703564388"
$ws.Cells.Item(1426, 5).Value = 1190

# Row 1427
$ws.Cells.Item(1427, 1).Value = 44448
$ws.Cells.Item(1427, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1427, 2).Value = "'5101004"
$ws.Cells.Item(1427, 3).Value = 3011
$ws.Cells.Item(1427, 4).Value = "Order 5101004 Swish +46704483544"
$ws.Cells.Item(1427, 6).Value = 1176.79

# Row 1428
$ws.Cells.Item(1428, 1).Value = 44448
$ws.Cells.Item(1428, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1428, 2).Value = "'5101004"
$ws.Cells.Item(1428, 3).Value = 2611
$ws.Cells.Item(1428, 4).Value = "Order 5101004 Swish +46704483544"
$ws.Cells.Item(1428, 6).Value = 141.21

# Row 1429
$ws.Cells.Item(1429, 1).Value = 44448
$ws.Cells.Item(1429, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1429, 2).Value = "'5101004"
$ws.Cells.Item(1429, 3).Value = 1930
$ws.Cells.Item(1429, 4).Value = "Order 5101004 Swish +46704483544"
$ws.Cells.Item(1429, 5).Value = 1318

# Row 1430
$ws.Cells.Item(1430, 1).Value = 44450
$ws.Cells.Item(1430, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1430, 3).Value = 4010
$ws.Cells.Item(1430, 4).Value = "NGROCERIES K0135"
$ws.Cells.Item(1430, 5).Value = 546.43

# Row 1431
$ws.Cells.Item(1431, 1).Value = 44450
$ws.Cells.Item(1431, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1431, 3).Value = 2645
$ws.Cells.Item(1431, 4).Value = "NGROCERIES K0135"
$ws.Cells.Item(1431, 5).Value = 65.57

# Row 1432
$ws.Cells.Item(1432, 1).Value = 44450
$ws.Cells.Item(1432, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1432, 3).Value = 1930
$ws.Cells.Item(1432, 4).Value = "NGROCERIES K0135"
$ws.Cells.Item(1432, 6).Value = 612

# Row 1433
$ws.Cells.Item(1433, 1).Value = 44450
$ws.Cells.Item(1433, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1433, 3).Value = 6400
$ws.Cells.Item(1433, 4).Value = "FACEBK 48MLK53Z62 K6885"
$ws.Cells.Item(1433, 5).Value = 200

# Row 1434
$ws.Cells.Item(1434, 1).Value = 44450
$ws.Cells.Item(1434, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1434, 4).Value = "FACEBK 48MLK53Z62 K6885"
$ws.Cells.Item(1434, 5).Value = 0

# Row 1435
$ws.Cells.Item(1435, 1).Value = 44450
$ws.Cells.Item(1435, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1435, 3).Value = 1930
$ws.Cells.Item(1435, 4).Value = "FACEBK 48MLK53Z62 K6885"
$ws.Cells.Item(1435, 6).Value = 200

# Row 1436
$ws.Cells.Item(1436, 1).Value = 44451
$ws.Cells.Item(1436, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1436, 2).Value = "'6121637"
$ws.Cells.Item(1436, 3).Value = 3011
$ws.Cells.Item(1436, 4).Value = "Order 6121637 Swish +46735689616"
$ws.Cells.Item(1436, 6).Value = 1062.5

# Row 1437
$ws.Cells.Item(1437, 1).Value = 44451
$ws.Cells.Item(1437, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1437, 2).Value = "'6121637"
$ws.Cells.Item(1437, 3).Value = 2611
$ws.Cells.Item(1437, 4).Value = "Order 6121637 Swish +46735689616"
$ws.Cells.Item(1437, 6).Value = 127.5

# Row 1438
$ws.Cells.Item(1438, 1).Value = 44451
$ws.Cells.Item(1438, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1438, 2).Value = "'6121637"
$ws.Cells.Item(1438, 3).Value = 1930
$ws.Cells.Item(1438, 4).Value = "Order 6121637 Swish +46735689616"
$ws.Cells.Item(1438, 5).Value = 1190

# Row 1439
$ws.Cells.Item(1439, 1).Value = 44451
$ws.Cells.Item(1439, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1439, 3).Value = 4010
$ws.Cells.Item(1439, 4).Value = "SNABBGROSS SOLNA K0135"
$ws.Cells.Item(1439, 5).Value = 1516.12

# Row 1440
$ws.Cells.Item(1440, 1).Value = 44451
$ws.Cells.Item(1440, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1440, 3).Value = 2645
$ws.Cells.Item(1440, 4).Value = "SNABBGROSS SOLNA K0135"
$ws.Cells.Item(1440, 5).Value = 181.93

# Row 1441
$ws.Cells.Item(1441, 1).Value = 44451
$ws.Cells.Item(1441, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1441, 3).Value = 1930
$ws.Cells.Item(1441, 4).Value = "SNABBGROSS SOLNA K0135"
$ws.Cells.Item(1441, 6).Value = 1698.05

# Row 1442
$ws.Cells.Item(1442, 1).Value = 44451
$ws.Cells.Item(1442, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1442, 3).Value = 4010
$ws.Cells.Item(1442, 4).Value = "SRI LANKA LIVS K0135"
$ws.Cells.Item(1442, 5).Value = 705.18

# Row 1443
$ws.Cells.Item(1443, 1).Value = 44451
$ws.Cells.Item(1443, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1443, 3).Value = 2645
$ws.Cells.Item(1443, 4).Value = "SRI LANKA LIVS K0135"
$ws.Cells.Item(1443, 5).Value = 84.62

# Row 1444
$ws.Cells.Item(1444, 1).Value = 44451
$ws.Cells.Item(1444, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1444, 3).Value = 1930
$ws.Cells.Item(1444, 4).Value = "SRI LANKA LIVS K0135"
$ws.Cells.Item(1444, 6).Value = 789.8
